$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Galicia (row 8) figures
$ws.Range("B8").Value = 8299
$ws.Range("C8").Value = 1548
$ws.Range("D8").Value = 6399
$ws.Range("E8").Value = 352

# 2. Update Navarra (row 12) figures
$ws.Range("B12").Value = 4697
$ws.Range("C12").Value = 1123
$ws.Range("D12").Value = 3189
$ws.Range("E12").Value = 385

# 3. Insert a new row at 19 (pushes existing rows 19-67 down to 20-68)
$ws.Rows(19).Insert()

# 4. Fill in the new row with the Extremadura data
$ws.Range("A19").Value = "Extremadura"
$ws.Range("B19").Value = 3186
$ws.Range("C19").Value = 858
$ws.Range("D19").Value = 1939
$ws.Range("E19").Value = 389

# 5. Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 15:52"
